$d = $word.ActiveDocument

# wdReplaceOne = 1, wdFindContinue = 1 ... constants used literally below:
#   Find.Execute(FindText, MatchCase, MatchWholeWord, MatchWildcards,
#                MatchSoundsLike, MatchAllWordForms, Forward, Wrap,
#                Format, ReplaceWith, Replace)
#   Wrap = 1 (wdFindContinue), Replace = 2 (wdReplaceAll)

# 1) Fusili quantity: 1500g -> 2000g
$d.Content.Find.Execute("1500g de ", $false, $false, $false, $false, $false, `
  $true, 1, $false, "2000g de ", 2) | Out-Null

# 2) Onions: 800g blanc -> 1000g rouges
$d.Content.Find.Execute("800g  d oignons blanc en lamelles", $false, $false, $false, $false, $false, `
  $true, 1, $false, "1000g  d oignons rouges en lamelles", 2) | Out-Null

# 3) Yellow peppers: 500g -> 700g
$d.Content.Find.Execute("500g de poivrons jaunes en lamelles", $false, $false, $false, $false, $false, `
  $true, 1, $false, "700g de poivrons jaunes en lamelles", 2) | Out-Null

# 4) Red peppers: 500g -> 700g
$d.Content.Find.Execute("500g de poivrons rouges en lamelles", $false, $false, $false, $false, $false, `
  $true, 1, $false, "700g de poivrons rouges en lamelles", 2) | Out-Null

# 5) Tomatoes: 800g -> 1000g
$d.Content.Find.Execute("800g de tomates en dés", $false, $false, $false, $false, $false, `
  $true, 1, $false, "1000g de tomates en dés", 2) | Out-Null

# 6) Jalapenos quantity: 65g -> 100g
$d.Content.Find.Execute("65g de ", $false, $false, $false, $false, $false, `
  $true, 1, $false, "100g de ", 2) | Out-Null

# 7) Garlic: "1 tète d ail  hachée" -> "2  tète d ail  hachée 20gousses environ"
$d.Content.Find.Execute("1 tète d ail  hachée", $false, $false, $false, $false, $false, `
  $true, 1, $false, "2  tète d ail  hachée 20gousses environ", 2) | Out-Null

# 8) Coulis de tomates: append dilution note
$d.Content.Find.Execute("680ml de coulis de tomates", $false, $false, $false, $false, $false, `
  $true, 1, $false, "680ml de coulis de tomates 369 ml pate tomate + 369 ml d eau", 2) | Out-Null

# 9) Cooking time for peppers: 5 -> 10 minutes
$d.Content.Find.Execute("s poivrons pendant environ 5 minutes.", $false, $false, $false, $false, $false, `
  $true, 1, $false, "s poivrons pendant environ 10 minutes.", 2) | Out-Null

# 10) Garlic cooking time: "1 à 2 minutes" -> "4a5 minutes"
$d.Content.Find.Execute(" la cuisson encore 1 à 2 minutes", $false, $false, $false, $false, $false, `
  $true, 1, $false, " la cuisson encore 4a5 minutes", 2) | Out-Null

# 11) Drop the stray Spanish-language override on the "jalapenos" run in the
#     instructions paragraph (it should just inherit the surrounding fr-CA
#     formatting). Locate the *second* occurrence of "jalapenos" (the first
#     one is in the ingredients list).
$first = $d.Content
$first.Find.Execute("jalapenos", $false, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$second = $d.Range($first.End, $d.Content.End)
$second.Find.Execute("jalapenos", $false, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$second.LanguageID = "fr-CA"
